# Add a new parameter row ("init_pr_ever_self_harmed_if_ever_depr") to the
# parameter_values sheet, inserted as row 15 (pushing base_3m_prob_depr and
# everything below it down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 15; existing rows 15-35 shift down to 16-36,
# keeping their own values/styles/row-heights intact.
$ws.Rows("15").Insert() | Out-Null

# The freshly-inserted row has no formatting of its own yet - copy the
# look of the rows immediately above it (same section of the table,
# style index 4 for the parameter name cell, 5 for the value cell, and
# D column with no explicit style) onto the new row 15.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("D2").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Fill in the new parameter's data.
$ws.Range("A15").Value = "init_pr_ever_self_harmed_if_ever_depr"
$ws.Range("B15").Value = 0.004
$ws.Range("D15").Value = "consistent with rate of incident self harm "

# Match the saved selection/active cell from the authored workbook.
$ws.Range("A16").Select() | Out-Null
